$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F, matching formatting of existing header row
$ws.Range("F1").Value = "chr"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# chr values for rows 2-13
$chrValues = @{
    2  = 10
    3  = 9
    4  = 6
    5  = 6
    6  = 8
    7  = 6
    8  = 6
    9  = 6
    10 = 6
    11 = 17
    12 = 6
    13 = 14
}

foreach ($row in $chrValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $chrValues[$row]
}
